$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, copying the format from the neighboring header cell (G1)
# so it picks up the same bold/border/centered style used by the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column H2:H4 ("Save" indicator values)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
